# Rework the "calibration" sheet's parameter table:
#  - E4 becomes a reference to "d52" (string) instead of the stray numeric 3
#  - rows 5-7 gain a "previous d-value" (col D) / "next d-value" (col E) pair
# and move the active selection to E9 (next empty row under the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("calibration")

$ws.Range("E4").Value = "d52"

$ws.Range("D5").Value = "a1"
$ws.Range("E5").Value = "f3"

$ws.Range("D6").Value = "f3"
$ws.Range("E6").Value = "f3"

$ws.Range("D7").Value = "g6"
$ws.Range("E7").Value = "d2"

[void]$ws.Range("E9").Select()
